$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "dot"
$ws.Range("B12").Value = "v.v"

$ws.Range("A13").Value = "comma"
$ws.Range("B13").Value = "v,v"

$ws.Range("A14").Value = "trailing1"
$ws.Range("A15").Value = "trailing2"
$ws.Range("B15").Value = "vv2  "
$ws.Range("B14").Value = "vv1 "

$ws.Range("H20").Select()
